$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header/style template cell so new cells inherit style index 1 (bold, border, centered)
$ws.Range("B1").Copy()
$ws.Range("B1:T1").PasteSpecial(-4122)

# Header row (row 1)
$ws.Range("B1").Value = "OBSV_PreC_Location_G_other"
$ws.Range("C1").Value = "OBSV_same_space_count_D_10p"
$ws.Range("D1").Value = "OBSV_PreC_meetingsize0ffice_D_10plus"
$ws.Range("E1").Value = "OBSV_Now_Interact_Work_within_org"
$ws.Range("F1").Value = "OBSV_PreC_interative_B_remote"
$ws.Range("G1").Value = "Collaborative Work Enivornmental Preferences"
$ws.Range("H1").Value = "Current Focused Work "
$ws.Range("I1").Value = "Focused Work Environmental Preferences"
$ws.Range("J1").Value = "Office Satisfaction"
$ws.Range("K1").Value = "Pre Covid In Person Collab Work (1-3 people)"
$ws.Range("L1").Value = "Pre Covid In Person Collab Work (4-10+ people)"
$ws.Range("M1").Value = "Pre-Covid Hybrid Meetings"
$ws.Range("N1").Value = "Pre-Covid In Person Meetings"
$ws.Range("O1").Value = "Pre-Covid Remote Meetings"
$ws.Range("P1").Value = "Future Workplace Features Preference"
$ws.Range("Q1").Value = "Environment Productivity Impact"
$ws.Range("R1").Value = "Workplace Tech Features Importance"
$ws.Range("S1").Value = "Workplace Preference"
$ws.Range("T1").Value = "Office Workspace Preference"

# Data rows 2-9
# row 2
$ws.Range("B2").Value = 713
$ws.Range("C2").Value = 713
$ws.Range("D2").Value = 713
$ws.Range("E2").Value = 710
$ws.Range("F2").Value = 713
$ws.Range("G2").Value = 713
$ws.Range("H2").Value = 713
$ws.Range("I2").Value = 713
$ws.Range("J2").Value = 713
$ws.Range("K2").Value = 713
$ws.Range("L2").Value = 713
$ws.Range("M2").Value = 713
$ws.Range("N2").Value = 713
$ws.Range("O2").Value = 713
$ws.Range("P2").Value = 713
$ws.Range("Q2").Value = 713
$ws.Range("R2").Value = 713
$ws.Range("S2").Value = 713
$ws.Range("T2").Value = 713
# row 3
$ws.Range("B3").Value = 0.6719200561009818
$ws.Range("C3").Value = 0.7095826426962472
$ws.Range("D3").Value = 0.2044419955920657
$ws.Range("E3").Value = 0.1501680322373501
$ws.Range("F3").Value = 0.8729593267882189
$ws.Range("G3").Value = 5.349021875292828
$ws.Range("H3").Value = 0.8504638108155419
$ws.Range("I3").Value = 4.084950709134868
$ws.Range("J3").Value = 4.122002820874472
$ws.Range("K3").Value = 0.2995980052993433
$ws.Range("L3").Value = 0.3589675994399235
$ws.Range("M3").Value = 0.358924131181476
$ws.Range("N3").Value = 0.4719637261585272
$ws.Range("O3").Value = 0.2856462039980057
$ws.Range("P3").Value = 2.202516432955356
$ws.Range("Q3").Value = 2.471004243281471
$ws.Range("R3").Value = 4.334273624823695
$ws.Range("S3").Value = 0.3784001611928269
$ws.Range("T3").Value = 0.1876065794848072
# row 4
$ws.Range("B4").Value = 0.4689307590702963
$ws.Range("C4").Value = 0.4251315873051083
$ws.Range("D4").Value = 0.374351734506099
$ws.Range("E4").Value = 0.2503496512593084
$ws.Range("F4").Value = 0.2831996551757824
$ws.Range("G4").Value = 1.50337978108837
$ws.Range("H4").Value = 0.2500109680649201
$ws.Range("I4").Value = 1.498473441905509
$ws.Range("J4").Value = 0.8013303159577678
$ws.Range("K4").Value = 0.1428528832329196
$ws.Range("L4").Value = 0.2793626181254321
$ws.Range("M4").Value = 0.2707302378835066
$ws.Range("N4").Value = 0.2974585252053411
$ws.Range("O4").Value = 0.3646232645609411
$ws.Range("P4").Value = 0.2965698816144888
$ws.Range("Q4").Value = 1.166024454252978
$ws.Range("R4").Value = 0.7669225226722741
$ws.Range("S4").Value = 0.247293969853232
$ws.Range("T4").Value = 0.05312349050094685
# row 5
$ws.Range("B5").Value = [double]"-2.220446049250313e-16"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = [double]"-2.220446049250313e-16"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 1
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0.1666666666666666
# row 6
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0.15
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.482344266190199
$ws.Range("H6").Value = 0.75
$ws.Range("I6").Value = 3.142857142857143
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 0.2222222222222222
$ws.Range("L6").Value = 0.1176470588235294
$ws.Range("M6").Value = 0.15
$ws.Range("N6").Value = 0.2
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 1.939956595129009
$ws.Range("Q6").Value = 1.666666666666667
$ws.Range("R6").Value = 4
$ws.Range("S6").Value = 0.1428571428571428
$ws.Range("T6").Value = 0.1666666666666667
# row 7
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.428571428571429
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 4.084950709134867
$ws.Range("J7").Value = 4
$ws.Range("K7").Value = 0.2995980052993432
$ws.Range("L7").Value = 0.3589675994399235
$ws.Range("M7").Value = 0.35
$ws.Range("N7").Value = 0.4719637261585272
$ws.Range("O7").Value = 0.1
$ws.Range("P7").Value = 2.164324610752314
$ws.Range("Q7").Value = 2.5
$ws.Range("R7").Value = 4.5
$ws.Range("S7").Value = 0.2857142857142857
$ws.Range("T7").Value = 0.1666666666666667
# row 8
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0.15
$ws.Range("E8").Value = 0.25
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.285714285714286
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 0.4
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 0.5
$ws.Range("N8").Value = 0.7000000000000001
$ws.Range("O8").Value = 0.4
$ws.Range("P8").Value = 2.384615384615385
$ws.Range("Q8").Value = 3.333333333333333
$ws.Range("R8").Value = 5
$ws.Range("S8").Value = 0.5714285714285714
$ws.Range("T8").Value = 0.1666666666666667
# row 9
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 5
$ws.Range("K9").Value = 0.5
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 2.846153846153846
$ws.Range("Q9").Value = 5
$ws.Range("R9").Value = 5
$ws.Range("S9").Value = 1
$ws.Range("T9").Value = 0.3221891441597311
